# Aniadido la tabla excel
# Adds one new inventory row (row 5) to the "Inventario" sheet/table:
#   E23 | 119 | 134-6479 | Digilent mod LVLSHFT Logic Level Shifter | 5 | PHR 22-T 01 | 1
# mirroring the layout/formatting already used by rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row values -------------------------------------------------
$ws.Range("A5").Value = "E23"
$ws.Range("B5").Value = 119
$ws.Range("C5").Value = "134-6479"
$ws.Range("D5").Value = "Digilent mod LVLSHFT Logic Level Shifter"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "PHR 22-T 01"
$ws.Range("G5").Value = 1

# --- Formatting: match F2:F4 / G2 centered style ---------------------
# -4108 = xlCenter
$ws.Range("F5").HorizontalAlignment = -4108
$ws.Range("G5").HorizontalAlignment = -4108

# --- Hyperlink on the RS code, same pattern as C3 / C4 ---------------
$ws.Hyperlinks.Add($ws.Range("C5"), "https://es.rs-online.com/web/p/convertidores-de-nivel-logico/1346479")
$ws.Range("C5").Style = "Hipervínculo"

# --- Keep the active selection on the newly added cell ---------------
[void]$ws.Range("G5").Select()
